$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force B10/C10 to be treated as text (not auto-converted to a date/number)
# while typing the values, then drop the number-format override so the
# cells end up with no explicit style - matching the other data rows.
$ws.Range("B10:C10").NumberFormat = "@"

$ws.Range("A10").Value = 1582848000
$ws.Range("B10").Value = "2020-02-28"
$ws.Range("C10").Value = "03030"
$ws.Range("D10").Value = "IDBTECH"
$ws.Range("E10").Value = 0.1
$ws.Range("F10").Value = 0.1
$ws.Range("G10").Value = 0.1
$ws.Range("H10").Value = 0.1
$ws.Range("I10").Value = "-"

$ws.Range("B10:C10").ClearFormats()
